# Weekly price-sheet update: a new weekly record for Ciboulette (Vega
# Central Mapocho de Santiago) is inserted as row 553, pushing the
# existing rows 553:600 down to 554:601 (dimension grows to A1:R601).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 553, shifting 553:600 -> 554:601.
$ws.Rows("553").Insert()

# Populate the newly inserted row 553 with this week's record.
$ws.Range("A553").Value = 9
$ws.Range("B553").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C553").Value = "Metropolitana"
$ws.Range("D553").Value = 45132
$ws.Range("E553").Value = 13
$ws.Range("F553").Value = 100112039
$ws.Range("G553").Value = "Ciboulette"
$ws.Range("H553").Value = "Sin especificar"
$ws.Range("I553").Value = "Primera"
$ws.Range("J553").Value = 340
$ws.Range("K553").Value = 1500
$ws.Range("L553").Value = 1600
$ws.Range("M553").Value = 1550
$ws.Range("N553").Value = "$/docena de atados"
$ws.Range("O553").Value = "Región Metropolitana"
$ws.Range("P553").Value = 517
$ws.Range("Q553").Value = 3
$ws.Range("R553").Value = "Hortaliza"
